$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "w1"
$ws.Range("H9").Value = "assignments/#draft; homework-2"
$ws.Range("I9").Value = "Draft; HW2"
$ws.Range("D3").Select()
